$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").ClearContents()
$ws.Range("C3").ClearContents()

$ws.Range("E5").Value = 4.10977504614245

$ws.Range("C8").Value = -1.479696720105184
$ws.Range("E8").Value = 2.503951807923088

$ws.Range("E9").Value = 1.194160460927862

$ws.Range("E10").Value = 1.459149667419735

$ws.Range("C11").Value = 2.192778679161966

$ws.Range("C12").Value = 3.408364488606752
$ws.Range("E12").Value = 2.570658574505447

$ws.Range("E13").Value = 2.479713128614192

$ws.Range("C15").Value = 1.666553973046025
$ws.Range("E15").Value = 0.7837484735427891

$ws.Range("C16").Value = 1.879266440112781
$ws.Range("E16").Value = 1.897396692213427

$ws.Range("C17").Value = -2.620683231370935

$ws.Range("C18").Value = -3.036556262700263

$ws.Range("E19").Value = -0.370811510370217
